$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

function Get-ParagraphByText {
    param($doc, [string]$needle)
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    throw "Paragraph containing '$needle' not found"
}

# --- Change 1: split "AUC2  (not sure...)" run, wrap "2  (" in gramStart/gramEnd proofErr ---
$p1xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="492CAACB" w14:textId="663C6EB3" w:rsidR="000D1DBC" w:rsidRPr="00D75C19" w:rsidRDefault="000D1DBC" w:rsidP="000D1DBC"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00D75C19"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Column </w:t></w:r><w:r w:rsidRPr="00D75C19"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>12</w:t></w:r><w:r w:rsidRPr="00D75C19"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="00D75C19" w:rsidRPr="00D75C19"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>AUC</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2  (</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>not sure what’s difference between both AUC me</w:t></w:r><w:r w:rsidR="00D75C19"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>thods</w:t></w:r><w:r w:rsidR="00D75C19" w:rsidRPr="00D75C19"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r></w:p>
'@

$para1 = Get-ParagraphByText $d "AUC2  (not sure"
$n1 = $para1.Range.InsertXML($p1xml)
Write-Output "change1 InsertXML n=$n1"

# --- Change 2: split "NOTE: The amount of ripples..." run, wrap "amount" in gramStart/gramEnd proofErr ---
$p2xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0C8053D4" w14:textId="319789D5" w:rsidR="006937F6" w:rsidRDefault="006937F6" w:rsidP="00EB6E4E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006937F6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Ripple_Waveforms_Rat_OS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006937F6"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: T</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">hese files contain the </w:t></w:r><w:r w:rsidR="0043655B"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">raw </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">ripple waveforms with the actual duration of each ripple. </w:t></w:r><w:r w:rsidR="000523D5"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">NOTE: The </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>amount</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> of ripples found here may differ from that of the GC files. This is because the traces extracted for the GC file are 6-second long, so in case a ripple is too close to a sleep stage transition and there are no 6 seconds of signal available, this ripple is then discarded from the GC files.</w:t></w:r><w:r w:rsidR="00B90D39"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> The structure is 1x9 and for each trial one would find X columns representing X ripples. For each ripples the signal length is Dx1 with D being the individual duration of each ripple. </w:t></w:r><w:r w:rsidR="009C6E79"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> The data is not filtered.</w:t></w:r></w:p>
'@

$para2 = Get-ParagraphByText $d "NOTE: The amount of ripples"
$n2 = $para2.Range.InsertXML($p2xml)
Write-Output "change2 InsertXML n=$n2"

# --- Change 3: bold the "T_cell" bullet and append new bold sentence ---
$p3xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="78BA74C0" w14:textId="662B7B60" w:rsidR="00565446" w:rsidRDefault="00565446" w:rsidP="00902416"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>T_cell</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: Contains the same data as T but instead of being a table it is a cell array. </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>This is currently the variable needed in the python script umap_test.py</w:t></w:r></w:p>
'@

$para3 = Get-ParagraphByText $d "T_cell: Contains the same data as T"
$n3 = $para3.Range.InsertXML($p3xml)
Write-Output "change3 InsertXML n=$n3"

Write-Output "done"
